$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All values below are text-formatted in the source sheet (inline strings),
# so they are written as strings (prefixed where needed) to avoid Excel
# auto-converting them to numeric/percentage types.

$updates = @{
    "D2"  = "301.13"
    "E2"  = "-1.28%"
    "D3"  = "31.36"
    "E3"  = "-2.95%"
    "D4"  = "5.138"
    "E4"  = "-3.22%"
    "D5"  = "0.07395"
    "E5"  = "-2.91%"
    "D6"  = "2.238"
    "E6"  = "25.07%"
    "D7"  = "7.930"
    "E7"  = "0.48%"
    "D8"  = "3.825"
    "E8"  = "-1.31%"
    "D9"  = "0.9195"
    "E9"  = "-1.24%"
    "D10" = "0.1711"
    "E10" = "-0.13%"
    "D11" = "0.07564"
    "E11" = "-5.21%"
    "D12" = "0.08079"
    "E12" = "0.46%"
    "D13" = "0.02994"
    "E13" = "-1.92%"
    "D14" = "0.09920"
    "E14" = "-0.14%"
    "D15" = "0.001506"
    "E15" = "1.04%"
    "D16" = "0.006157"
    "E16" = "-2.34%"
    "D17" = "3.474"
    "E17" = "0.87%"
    "D18" = "2.227"
    "E18" = "-0.21%"
    "D19" = "0.3261"
    "E19" = "-1.21%"
    "D20" = "0.1318"
    "E20" = "-1.96%"
    "D21" = "4.652"
    "E21" = "2.30%"
    "D22" = "0.04635"
    "E22" = "0.75%"
    "D23" = "0.1566"
    "E23" = "-3.08%"
    "D25" = "0.004482"
    "E25" = "-0.41%"
    "D26" = "0.0001299"
    "E26" = "-6.91%"
    "D27" = "0.0003426"
    "E27" = "92.53%"
    "D39" = "0.01732"
    "E39" = "-0.34%"
    "D40" = "0.04497"
    "E40" = "-1.07%"
    "D41" = "0.007300"
    "E41" = "4.78%"
    "D42" = "0.1349"
    "E42" = "-0.96%"
    "D43" = "0.002228"
    "E43" = "8.01%"
    "E44" = "-23.61%"
    "D45" = "0.00006260"
    "E45" = "1.98%"
    "B46" = "BOLO"
    "C46" = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
    "D46" = "0.8165"
    "E46" = "15.15%"
    "B47" = "CoinbaseStockToken"
    "C47" = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
    "D47" = "0.009992"
    "E47" = "-18.10%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
